$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Career Projection")

# New columns' header cells first (G, H, I) so the new shared strings are
# interned in the same order the committed workbook used, then the
# E1/F1 pair (Move Month / Cost of Living).
$ws.Range("G1").Value = "State Taxes"
$ws.Range("H1").Value = "Dependents"
$ws.Range("I1").Value = "Married"
$ws.Range("E1").Value = "Move Month"
$ws.Range("F1").Value = "Cost of Living"

# First data row's year moved from 2022 to 2020 (a "boring" pre-career year).
$ws.Range("A2").Value = 2020

# New, narrower helper columns for the added fields.
$ws.Columns.Item(7).ColumnWidth = 10.17
$ws.Columns.Item(8).ColumnWidth = 9.5

# Selection moved onto the new E2 cell.
$ws.Activate()
$ws.Range("E2").Select()
